$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 11; $r++) {
    $formula = '="INSERT INTO "&A' + $r + '&" ([" & B$1 &"],["&C$1&"],["&D$1&"],["&E$1&"],["&F$1&"]) VALUES ( ''" & B' + $r + ' & "'',''" & C' + $r + ' & "'',''" & D' + $r + ' & "'' ,''" & E' + $r + ' & "'',''" & F' + $r + ' & "'')"'
    $ws.Range("G$r").Formula = $formula
}

$ws.Range("G14").Select()
